$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")
$ws.Activate()

# Item code was selected from the accessories list: F3 now holds the
# value previously in F4, and F4 is cleared out.
$ws.Range("F3").Value = "TB7SX6CC"
$ws.Range("F4").ClearContents()

# Move the active selection to H13, as left by the user after picking it.
$ws.Range("H13").Select()
